$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.264.04"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.358.39"
$ws.Range("E3").Value = "  +5.69%  "
$ws.Range("E4").Value = "  -0.61%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'232.88"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("E6").Value = "  -0.01%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").Value = "'68.02"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +9.19%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +2.84%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.0966"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +0.93%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'56.86"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.09%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").Value = "'26.37"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "2.709.59"
$ws.Range("E13").Value = "  +5.62%  "
$ws.Range("E14").Value = "  -0.65%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'15.75"
$ws.Range("D15").Style = $origStyle
$ws.Range("E16").Value = "  +3.20%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").Value = "'0.844"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "2.357.74"
$ws.Range("E18").Value = "  +5.13%  "
$ws.Range("D19").Value = "43.342.01"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("E20").Value = "  +0.14%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").Value = "'74.03"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +2.23%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'6.28"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +4.92%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'249.11"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +0.64%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'3.97"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +17.21%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("E28").Value = "  +1.49%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'22.44"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +8.44%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").Value = "'172.42"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +1.35%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'1.55"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +10.53%  "
$ws.Range("E32").Value = "  -7.42%  "
$ws.Range("E33").Value = "  +0.19%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").Value = "'5.03"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +6.58%  "
$ws.Range("E35").Value = "  -0.01%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").Value = "'5.06"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +3.99%  "
$ws.Range("E37").Value = "  +11.02%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").Value = "'6.52"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  -0.10%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.0255"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("E41").Value = "  +9.61%  "
$ws.Range("E42").Value = "  -0.07%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = "'18.30"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +8.08%  "
$ws.Range("E44").Value = "  +10.28%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").Value = "'98.97"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").Value = "'1.21"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +2.89%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'4.43"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +2.17%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.0954"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "1.450.67"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "2.581.06"
$ws.Range("E50").Value = "  +5.76%  "
$ws.Range("E51").Value = "  -1.05%  "
